$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D accepts numeric-looking values (e.g. "1.002") as plain text,
# matching the source data which stores these as inline strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Updated coin prices (column D)
$ws.Range("D2").Value = "26.038.61"
$ws.Range("D3").Value = "1.716.22"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "330.24"
$ws.Range("D6").Value = "0.9979"
$ws.Range("D7").Value = "0.3690"
$ws.Range("D8").Value = "49.88"
$ws.Range("D9").Value = "0.3324"
$ws.Range("D10").Value = "1.181"
$ws.Range("D11").Value = "0.07475"
$ws.Range("D12").Value = "0.9981"
$ws.Range("D14").Value = "20.02"
$ws.Range("D15").Value = "6.910"
$ws.Range("D16").Value = "1.713.97"
$ws.Range("D17").Value = "0.00001076"
$ws.Range("D18").Value = "0.06631"
$ws.Range("D19").Value = "81.98"
$ws.Range("D20").Value = "0.9980"
$ws.Range("D21").Value = "16.32"
$ws.Range("D22").Value = "6.060"
$ws.Range("D23").Value = "12.98"
$ws.Range("D24").Value = "25.992.51"
$ws.Range("D25").Value = "2.477"
$ws.Range("D26").Value = "2.481"
$ws.Range("D27").Value = "149.83"
$ws.Range("D28").Value = "19.26"
$ws.Range("D30").Value = "1.904.94"
$ws.Range("D31").Value = "128.99"
$ws.Range("D32").Value = "4.108"
$ws.Range("D33").Value = "5.955"
$ws.Range("D34").Value = "0.08527"
$ws.Range("D35").Value = "1.714"
$ws.Range("D36").Value = "12.90"
$ws.Range("D37").Value = "5.345"
$ws.Range("D38").Value = "1.281"
$ws.Range("D39").Value = "0.06190"
$ws.Range("D40").Value = "0.02282"
$ws.Range("D41").Value = "0.2128"
$ws.Range("D42").Value = "8.497"
$ws.Range("D43").Value = "14.43"
$ws.Range("D44").Value = "0.6152"
$ws.Range("D45").Value = "0.9983"
$ws.Range("D46").Value = "3.833"
$ws.Range("D47").Value = "0.5863"
$ws.Range("D48").Value = "127.22"
$ws.Range("D49").Value = "2.005"
$ws.Range("D50").Value = "0.07258"
$ws.Range("D51").Value = "77.01"

# Restore default (Normal) cell style now that the text values are committed,
# so no stray number-format style lingers on these cells.
$ws.Range("D2:D51").Style = "Normal"

# Updated 1h volume percentages (column E)
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("E8").Value = "  +5.75%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +4.38%  "
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +3.28%  "
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("E24").Value = "  +5.59%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  +8.31%  "
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("E36").Value = "  +4.61%  "
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("E43").Value = "  +13.19%  "
$ws.Range("E44").Value = "  +3.67%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("E51").Value = "  +3.21%  "
